# Auto-generated Excel COM-interop edit script
# Applies per-cell numeric updates to the Golem_Profits-style leve profit sheets
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 773
$ws.Range("I6").Value = 792.5
$ws.Range("J6").Value = 695
$ws.Range("K6").Value = 2377.5
$ws.Range("L6").Value = 2085
$ws.Range("M6").Value = -2265.5
$ws.Range("N6").Value = -2309
$ws.Range("H12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()
$ws.Range("H116").Value = 3000
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 3000
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 3000
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -9884
$ws.Range("H125").Value = 817.5
$ws.Range("J125").Value = 548
$ws.Range("L125").Value = 4932
$ws.Range("N125").Value = -9852

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 180
$ws.Range("I17").Value = 180
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 180
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -7
$ws.Range("N17").ClearContents()
$ws.Range("H32").Value = 1470.25
$ws.Range("I32").Value = 455.5
$ws.Range("K32").Value = 455.5
$ws.Range("M32").Value = -168.5
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("H50").Value = 17217.572
$ws.Range("I50").Value = 6848
$ws.Range("J50").Value = 24994.75
$ws.Range("K50").Value = 6848
$ws.Range("L50").Value = 24994.75
$ws.Range("M50").Value = -6134
$ws.Range("N50").Value = -26422.75
$ws.Range("H74").Value = 1270
$ws.Range("J74").Value = 1265
$ws.Range("L74").Value = 1265
$ws.Range("N74").Value = -3013
$ws.Range("H77").Value = 1270
$ws.Range("J77").Value = 1265
$ws.Range("L77").Value = 6325
$ws.Range("N77").Value = -15061
$ws.Range("H122").Value = 4285.875
$ws.Range("I122").Value = 4285.875
$ws.Range("K122").Value = 12857.625
$ws.Range("M122").Value = -10407.625
$ws.Range("H132").Value = 2000
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 6000
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -11060
$ws.Range("H135").Value = 3360666.8
$ws.Range("J135").Value = 3360666.8
$ws.Range("L135").Value = 3360666.8
$ws.Range("N135").Value = -3370806.8

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8292.05
$ws.Range("I31").Value = 2678.3125
$ws.Range("J31").Value = 30747
$ws.Range("K31").Value = 2678.3125
$ws.Range("L31").Value = 30747
$ws.Range("M31").Value = -2383.3125
$ws.Range("N31").Value = -31337
$ws.Range("H33").Value = 11746
$ws.Range("J33").Value = 40999.332
$ws.Range("L33").Value = 40999.332
$ws.Range("N33").Value = -41757.332
$ws.Range("H34").Value = 8292.05
$ws.Range("I34").Value = 2678.3125
$ws.Range("J34").Value = 30747
$ws.Range("K34").Value = 2678.3125
$ws.Range("L34").Value = 30747
$ws.Range("M34").Value = -2476.3125
$ws.Range("N34").Value = -31151
$ws.Range("H44").Value = 29998.6
$ws.Range("I44").Value = 29997
$ws.Range("J44").Value = 29999.666
$ws.Range("K44").Value = 29997
$ws.Range("L44").Value = 29999.666
$ws.Range("M44").Value = -29555
$ws.Range("N44").Value = -30883.666
$ws.Range("H55").Value = 9000
$ws.Range("I55").Value = 9000
$ws.Range("K55").Value = 9000
$ws.Range("M55").Value = -8685
$ws.Range("H56").Value = 13333.333
$ws.Range("I56").Value = 13333.333
$ws.Range("K56").Value = 13333.333
$ws.Range("M56").Value = -12488.333
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("N58").ClearContents()
$ws.Range("H60").Value = 31956.4
$ws.Range("J60").Value = 36672.5
$ws.Range("L60").Value = 36672.5
$ws.Range("N60").Value = -37694.5
$ws.Range("H94").Value = 3392.5
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 3392.5
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 3392.5
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -4294.5
$ws.Range("H107").Value = 927.1429
$ws.Range("J107").Value = 897
$ws.Range("L107").Value = 897
$ws.Range("N107").Value = -4737
$ws.Range("H132").Value = 2075.3333
$ws.Range("I132").Value = 1945.8182
$ws.Range("K132").Value = 5837.4546
$ws.Range("M132").Value = -3307.4546
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2051.5
$ws.Range("J68").Value = 2051.5
$ws.Range("L68").Value = 6154.5
$ws.Range("N68").Value = -7776.5
$ws.Range("H71").Value = 2051.5
$ws.Range("J71").Value = 2051.5
$ws.Range("L71").Value = 18463.5
$ws.Range("N71").Value = -26575.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 3328677
$ws.Range("I14").Value = 3933072.8
$ws.Range("K14").Value = 3933072.8
$ws.Range("M14").Value = -3932904.8
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H102").Value = 67900.664
$ws.Range("I102").Value = 80480.8
$ws.Range("J102").Value = 5000
$ws.Range("K102").Value = 80480.8
$ws.Range("L102").Value = 5000
$ws.Range("M102").Value = -78858.8
$ws.Range("N102").Value = -8244
$ws.Range("H123").Value = 54875.75
$ws.Range("J123").Value = 54875.75
$ws.Range("L123").Value = 54875.75
$ws.Range("N123").Value = -59775.75
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 915.5
$ws.Range("I7").Value = 915.5
$ws.Range("K7").Value = 915.5
$ws.Range("M7").Value = -803.5
$ws.Range("H23").Value = 3000
$ws.Range("I23").Value = 3000
$ws.Range("K23").Value = 3000
$ws.Range("M23").Value = -2770
$ws.Range("H40").Value = 24183.334
$ws.Range("I40").Value = 27335.2
$ws.Range("J40").Value = 8424
$ws.Range("K40").Value = 27335.2
$ws.Range("L40").Value = 8424
$ws.Range("M40").Value = -27199.2
$ws.Range("N40").Value = -8696
$ws.Range("H45").Value = 15000
$ws.Range("I45").Value = 15000
$ws.Range("K45").Value = 15000
$ws.Range("M45").Value = -14593
$ws.Range("H55").Value = 306.15384
$ws.Range("I55").Value = 317.4
$ws.Range("K55").Value = 317.4
$ws.Range("M55").Value = -144.4
$ws.Range("H76").Value = 14288
$ws.Range("J76").Value = 14288
$ws.Range("L76").Value = 14288
$ws.Range("N76").Value = -14964
$ws.Range("H79").Value = 14288
$ws.Range("J79").Value = 14288
$ws.Range("L79").Value = 14288
$ws.Range("N79").Value = -16628
$ws.Range("H126").Value = 915.5
$ws.Range("I126").Value = 915.5
$ws.Range("K126").Value = 2746.5
$ws.Range("M126").Value = -276.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 475
$ws.Range("I81").Value = 475
$ws.Range("K81").Value = 950
$ws.Range("M81").Value = 111
$ws.Range("H84").Value = 475
$ws.Range("I84").Value = 475
$ws.Range("K84").Value = 4750
$ws.Range("M84").Value = 554
$ws.Range("H132").Value = 6000
$ws.Range("I132").Value = 6000
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 18000
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -15470
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 6067.8667
$ws.Range("I136").Value = 4534.5
$ws.Range("J136").Value = 7820.2856
$ws.Range("K136").Value = 13603.5
$ws.Range("L136").Value = 23460.8568
$ws.Range("M136").Value = -11053.5
$ws.Range("N136").Value = -28560.8568
